$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I18").Value = -0.1107587875764343
$ws.Range("J18").Value = 0.1116453926146825
$ws.Range("K18").Value = 0.1635099888121481
$ws.Range("L18").Value = 2.252950313514092
